$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: D15:I15 and M15:R15 switch from empty numeric cells to "-" inline strings
# with the right-aligned dash style (same style used by J15:L15 already).
$r1 = $ws.Range("D15:I15")
$r1.Value = "-"
$r1.HorizontalAlignment = -4152

$r2 = $ws.Range("M15:R15")
$r2.Value = "-"
$r2.HorizontalAlignment = -4152

# Row 70: J70:R70 switch from empty numeric cells to "-" inline strings.
$r3 = $ws.Range("J70:R70")
$r3.Value = "-"
$r3.HorizontalAlignment = -4152

# Row 73: J73:R73 switch from empty numeric cells to "-" inline strings.
$r4 = $ws.Range("J73:R73")
$r4.Value = "-"
$r4.HorizontalAlignment = -4152

# Row 76: J76:R76 switch from empty numeric cells to "-" inline strings.
$r5 = $ws.Range("J76:R76")
$r5.Value = "-"
$r5.HorizontalAlignment = -4152
